$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.133.38'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").Value = '3.321.85'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  +1.34%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '3.313.20'
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  +7.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.629'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000280'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '3.848.83'
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.334.79'
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '64.188.51'
$ws.Range("E19").Value = '  +1.54%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.983'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '452.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.89%  '
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.19%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '564.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '60.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.59%  '
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.140'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.366'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0730'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("D42").Value = '3.052.34'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.75'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.46%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0412'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.68%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.133'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.88%  '
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("E51").Value = '  +0.94%  '
